$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: "Emails" keyword step becomes "Compose"
$ws.Range("D4").Value = "Compose"

# Row 6: the "Compose" keyword step is cleared out (keyword column emptied)
$ws.Range("D6").ClearContents()

# The whole second test case block (TC_02, rows 7-11) is removed,
# leaving only the blank, still-styled A/B cells behind.
$ws.Range("A7:H11").ClearContents()

# Update page setup (paper size / orientation) as recorded for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moves to E4.
$ws.Range("E4").Select()
